# The sheet originally held 5 data rows (rows 2-5) of password/email test
# data.  The fix keeps only the first data row (row 2) and removes rows 3-5
# entirely (their cell values, styles and the shared strings they alone
# used all disappear together with them).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows("3:5").Delete()

# The remaining row's email cell is updated to a new value.
$ws.Range("D2").Value = "1asd@gmail.com23"

# Rows 3-5 carried hyperlinks (for D3:D5, D3, D4 and D5) that become stale
# once those rows are gone. Rebuild the hyperlink collection from scratch so
# only the one hyperlink belonging to D2 remains, pointing at the original
# mailto: address.
$ws.Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("D2"), "mailto:asd@gmail.com")
$ws.Range("D2").Style = "Hyperlink"

# Leave the selection on D3 (the cell below the last remaining row), matching
# where the user ended up after the edit.
$ws.Range("D3").Select()
